$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 170.88889
$ws.Range("I41").Value = 170.88889
$ws.Range("K41").Value = 170.88889
$ws.Range("M41").Value = 269.11111

$ws.Range("H53").Value = 403.9
$ws.Range("I53").Value = 503.42856
$ws.Range("K53").Value = 503.42856
$ws.Range("M53").Value = 133.57144

$ws.Range("H97").Value = 680
$ws.Range("J97").Value = 680
$ws.Range("L97").Value = 2040
$ws.Range("N97").Value = -3032

$ws.Range("H131").Value = 2061.6
$ws.Range("I131").Value = 2061.6
$ws.Range("K131").Value = 6184.799999999999
$ws.Range("M131").Value = -1144.799999999999

$ws.Range("H137").Value = 2771
$ws.Range("J137").Value = 3524.8333
$ws.Range("L137").Value = 10574.4999
$ws.Range("N137").Value = -15674.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1022.5
$ws.Range("I2").Value = 1022.5
$ws.Range("K2").Value = 1022.5
$ws.Range("M2").Value = -909.5

$ws.Range("H5").Value = 275.66666
$ws.Range("I5").Value = 233.1
$ws.Range("K5").Value = 233.1
$ws.Range("M5").Value = -121.1

$ws.Range("H6").Value = 12003567
$ws.Range("I6").Value = 10913091
$ws.Range("J6").Value = 15002375
$ws.Range("K6").Value = 10913091
$ws.Range("L6").Value = 15002375
$ws.Range("M6").Value = -10912918
$ws.Range("N6").Value = -15002721

$ws.Range("H61").Value = 2405.5386
$ws.Range("I61").Value = 2405.5386
$ws.Range("K61").Value = 2405.5386
$ws.Range("M61").Value = -2193.5386

$ws.Range("H97").Value = 949.4737
$ws.Range("I97").Value = 857.61536
$ws.Range("J97").Value = 1148.5
$ws.Range("K97").Value = 857.61536
$ws.Range("L97").Value = 1148.5
$ws.Range("M97").Value = -361.61536
$ws.Range("N97").Value = -2140.5

$ws.Range("H116").Value = 1022.5
$ws.Range("I116").Value = 1022.5
$ws.Range("K116").Value = 1022.5
$ws.Range("M116").Value = 1271.5

$ws.Range("H122").Value = 1420.0435
$ws.Range("I122").Value = 1224.8889
$ws.Range("K122").Value = 3674.6667
$ws.Range("M122").Value = -1224.6667

$ws.Range("H132").Value = 2974.9167
$ws.Range("I132").Value = 2877.889
$ws.Range("K132").Value = 8633.667000000001
$ws.Range("M132").Value = -6103.667000000001

$ws.Range("H136").Value = 2405.5386
$ws.Range("I136").Value = 2405.5386
$ws.Range("K136").Value = 7216.6158
$ws.Range("M136").Value = -4666.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1022.5
$ws.Range("I3").Value = 1022.5
$ws.Range("K3").Value = 1022.5
$ws.Range("M3").Value = -908.5

$ws.Range("H4").Value = 275.66666
$ws.Range("I4").Value = 233.1
$ws.Range("K4").Value = 233.1
$ws.Range("M4").Value = -118.1

$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H86").Value = 4333.6665
$ws.Range("I86").Value = 5000.5
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 5000.5
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -3877.5
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 4333.6665
$ws.Range("I89").Value = 5000.5
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 25002.5
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -19386.5
$ws.Range("N89").Value = -26232

$ws.Range("H94").Value = 424.7
$ws.Range("I94").Value = 424.7
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 424.7
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = 26.30000000000001

$ws.Range("H134").Value = 5772
$ws.Range("I134").Value = 1331.4546
$ws.Range("J134").Value = 30195
$ws.Range("K134").Value = 3994.3638
$ws.Range("L134").Value = 90585
$ws.Range("M134").Value = -1459.3638
$ws.Range("N134").Value = -95655

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 63.214287
$ws.Range("I7").Value = 53.555557
$ws.Range("J7").Value = 80.59999999999999
$ws.Range("K7").Value = 53.555557
$ws.Range("L7").Value = 80.59999999999999
$ws.Range("M7").Value = 59.444443
$ws.Range("N7").Value = -306.6

$ws.Range("H31").Value = 6677
$ws.Range("I31").Value = 4373
$ws.Range("J31").Value = 7103.6665
$ws.Range("K31").Value = 4373
$ws.Range("L31").Value = 7103.6665
$ws.Range("M31").Value = -4078
$ws.Range("N31").Value = -7693.6665

$ws.Range("H34").Value = 6677
$ws.Range("I34").Value = 4373
$ws.Range("J34").Value = 7103.6665
$ws.Range("K34").Value = 4373
$ws.Range("L34").Value = 7103.6665
$ws.Range("M34").Value = -4171
$ws.Range("N34").Value = -7507.6665

$ws.Range("H86").Value = 7000
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877

$ws.Range("H89").Value = 7000
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384

$ws.Range("H105").Value = 1340
$ws.Range("I105").Value = 1066.6666
$ws.Range("J105").Value = 1750
$ws.Range("K105").Value = 1066.6666
$ws.Range("L105").Value = 1750
$ws.Range("M105").Value = 680.3334
$ws.Range("N105").Value = -5244

$ws.Range("H132").Value = 1785
$ws.Range("I132").Value = 1383.5834
$ws.Range("K132").Value = 4150.7502
$ws.Range("M132").Value = -1620.7502

$ws.Range("H134").Value = 2132.75
$ws.Range("I134").Value = 763.2727
$ws.Range("K134").Value = 2289.8181
$ws.Range("M134").Value = 245.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 864.5714
$ws.Range("I139").Value = 864.5714
$ws.Range("K139").Value = 2593.7142
$ws.Range("M139").Value = 2546.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 22500
$ws.Range("J38").Value = 22500
$ws.Range("L38").Value = 22500
$ws.Range("N38").Value = -23426

$ws.Range("H80").Value = 4838.8
$ws.Range("I80").Value = 4798.5
$ws.Range("K80").Value = 4798.5
$ws.Range("M80").Value = -3800.5

$ws.Range("H83").Value = 4838.8
$ws.Range("I83").Value = 4798.5
$ws.Range("K83").Value = 23992.5
$ws.Range("M83").Value = -19000.5

$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 727.7273
$ws.Range("I46").Value = 550
$ws.Range("J46").Value = 875.8333
$ws.Range("K46").Value = 550
$ws.Range("L46").Value = 875.8333
$ws.Range("M46").Value = -362
$ws.Range("N46").Value = -1251.8333

$ws.Range("H132").Value = 3725.7273
$ws.Range("I132").Value = 4028.3
$ws.Range("K132").Value = 12084.9
$ws.Range("M132").Value = -9554.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1510
$ws.Range("I29").Value = 1510
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1510
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -1220

$ws.Range("H113").Value = 917.3333
$ws.Range("I113").Value = 773.1429000000001
$ws.Range("K113").Value = 2319.4287
$ws.Range("M113").Value = -149.4287000000004

$ws.Range("H126").Value = 3490.261
$ws.Range("I126").Value = 1417.8462
$ws.Range("K126").Value = 4253.5386
$ws.Range("M126").Value = -1783.5386
